# Apply the SmartRewards section rewrite (English -> updated copy, as part of the
# Vietnamese-translation commit). The original section was 4 paragraphs:
#   A) "SmartRewards are a price stabilization mechanism ... 15% block reward allocation."
#   B) "Here is a small example: if you buy SmartCash on the 24th ..."
#   C) "These payments will occur on the 25th of each month."
#   D) "All users need to move funds into addresses larger than 1000 SMART before the 25th..."
# It gets replaced by 8 paragraphs (4 with text, interleaved with 4 empty spacer
# paragraphs) carrying the new SmartRewards / Decentralized Distribution copy.

$d = $word.ActiveDocument

# Locate paragraph A (the long "SmartRewards ... " paragraph) by distinctive text.
$rngA = $d.Content
$foundA = $rngA.Find.Execute("Each wallet address holding at least 1000 SmartCash", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundA) {
    throw "Could not locate the original SmartRewards paragraph (A)."
}
$rngA.Expand(4) | Out-Null

# Locate paragraph D (the last "All users need to move funds..." paragraph) by distinctive text.
$rngD = $d.Content
$foundD = $rngD.Find.Execute("All users need to move funds into addresses larger than 1000 SMART", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundD) {
    throw "Could not locate the trailing 'All users need to move funds' paragraph (D)."
}
$rngD.Expand(4) | Out-Null

# Build one range spanning paragraphs A through D (inclusive) and replace their
# combined contents (pPr + runs) in a single shot via raw OOXML injection, so
# every run keeps its exact formatting (fonts/size/bold/italic/borders/etc.)
$full = $d.Range($rngA.Start, $rngD.End)

$xml = @'
<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>SmartRewards</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> are a price stabilization mechanism and a way to encourage long term holding. Long term holders are key to the project’s success since the SmartHive treasury needs SmartCash to appreciate in value in order to fund meaningful 3rd party proposals and help grow SmartCash into a successful global crypto-currency. Beginning at block 574,100, the current SmartRewards model will change to a Decentralized Distribution. SmartRewards will then be distributed after every </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">47,500 Blocks</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> to all eligible addresses. The payouts will begin 200 Blocks after the cycle ends and 1000 addresses will be paid every other block. SmartNodes will keep track on which addresses are eligible and which are not. You will earn SmartRewards on </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> address for which you hold the keys (web or desktop, SmartNode included!) which holds &gt;=1000 SMART for one month and does not make </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>any</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">outgoing transactions during that time. Please note, most exchanges do </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:i/><w:iCs/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> pay SmartRewards to their users, holding &gt;= 1000 SMART on an exchange does not guarantee a reward</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">The SmartRewards will come out of the 15% block reward allocation.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">Payment after every 47500 blocks starting at 574100. Typically, around the 25th of each month.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">All users need to move funds into addresses holding at least 1000 SMART before the snapshot to be counted.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">If you spend ANY amount from an address, it will be ineligible for SmartRewards until the next round.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="396" w:lineRule="atLeast"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="Times New Roman" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>
'@

$full.InsertXML($xml) | Out-Null
